$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Revert sheet name from template placeholder to plain Russian text
$ws.Name = "Контакты"

# Revert header cells from template placeholders to plain Russian text
$ws.Range("A1").Value = "Имя"
$ws.Range("B1").Value = "Адрес"
$ws.Range("C1").Value = "Квартира"
$ws.Range("D1").Value = "Телефон"
$ws.Range("E1").Value = "Почта"
